$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.370913
$ws.Range("H2").Value = 1.112739
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.002008666666667
$ws.Range("N2").Value = 12.006026
$ws.Range("O2").Value = 0.4834231243738785
$ws.Range("P2").Value = 0.4834231243738787
$ws.Range("Q2").Value = 1.484397040579333
$ws.Range("R2").Value = 13.359573365214
$ws.Range("S2").Value = 0.4834231243738785
$ws.Range("T2").Value = 0.4834231243738787

# Update row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.370913
$ws.Range("H3").Value = 1.112739
$ws.Range("M3").Value = 4.265473333333333
$ws.Range("N3").Value = 12.79642
$ws.Range("O3").Value = 0.5152483708764571
$ws.Range("P3").Value = 0.5152483708764573
$ws.Range("Q3").Value = 1.582119510486667
$ws.Range("R3").Value = 14.23907559438
$ws.Range("S3").Value = 0.5152483708764571
$ws.Range("T3").Value = 0.5152483708764573

# Add new row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna3"
$ws.Range("C4").Value = "Ephb1"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.370913
$ws.Range("H4").Value = 1.112739
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.010998
$ws.Range("N4").Value = 0.032994
$ws.Range("O4").Value = 0.001328504749664189
$ws.Range("P4").Value = 0.00132850474966419
$ws.Range("Q4").Value = 0.004079301174000001
$ws.Range("R4").Value = 0.036713710566
$ws.Range("S4").Value = 0.001328504749664189
$ws.Range("T4").Value = 0.00132850474966419
